# Css101.docx edit:
#   "Souradyuti paul" (single run, one paragraph)
# becomes two paragraphs:
#   "Souradyuti " + proofErr-wrapped "paul"
#   proofErr-wrapped "Iit" + " bhilai"   (new paragraph, inherits the _GoBack bookmark)
#
# w:proofErr (the wavy-underline "possible spelling error" markers) and the exact
# run-split aren't reachable through the normal Range/Find object model, so we
# build the target OOXML for each paragraph and splice it in with Range.InsertXML
# (the documented way to inject raw WordProcessingML through the COM OM).

$d = $word.ActiveDocument

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wordNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData>' +
           '</pkg:part></pkg:package>'
}

# --- Paragraph 1: split "Souradyuti paul" into "Souradyuti " + proofErr("paul") ---
$firstPara = $d.Paragraphs.First
$p1Range = $firstPara.Range.Duplicate
[void]$p1Range.MoveEnd(1, -1)   # wdCharacter, -1: drop the trailing paragraph mark

$p1Body = '<w:p>' +
            '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Souradyuti </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>paul</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
          '</w:p>'
[void]$p1Range.InsertXML((New-PkgXml $p1Body))

# InsertXML re-materializes any bookmark that used to sit inside the replaced
# range, so "_GoBack" (collapsed, right after "paul") is now back at the end of
# paragraph 1. It belongs at the end of the new paragraph 2 instead, so drop it
# here and re-create it there explicitly.
try {
    $goBack = $d.Bookmarks("_GoBack")
    [void]$goBack.Delete()
} catch {
    # no _GoBack bookmark present - nothing to relocate
}

# --- Paragraph 2 (new): proofErr("Iit") + " bhilai", plus the relocated bookmark ---
$p1 = $d.Paragraphs(1)
$insertAt = $p1.Range.End - 1   # position just before paragraph 1's paragraph mark
$p2Range = $d.Range($insertAt, $insertAt)

$p2Body = '<w:p>' +
            '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Iit</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> bhilai</w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
            '<w:bookmarkEnd w:id="0"/>' +
          '</w:p>'
[void]$p2Range.InsertXML((New-PkgXml $p2Body))
